$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    Appears on the Overview sheet (E2:E3 "zh-cn" column, F2:F3 "de-de"
#    column) and on the per-locale sheets' "Status" column (C2:C3).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback report generated: fill "Latest Target File" (I),
#    "Latest Handback File" (J) and "Latest Handback DateTime" (K) for
#    both data rows on the zh-cn and de-de sheets.
# ---------------------------------------------------------------------------

# zh-cn: target file "a.md" (hyperlinked), handback xliff + timestamp.
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b57f9966c4376dbfb47aec704b0c44192853980/e2e/a.md", "", "", "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b57f9966c4376dbfb47aec704b0c44192853980/e2e/a.md", "", "", "a.md") | Out-Null
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-29 18:40:31"
$wsZh.Range("K3").Value = "2016-08-29 18:40:31"

# de-de: target file "a.md" (hyperlinked), handback xliff + timestamp.
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b57f9966c4376dbfb47aec704b0c44192853980/e2e/a.md", "", "", "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b57f9966c4376dbfb47aec704b0c44192853980/e2e/a.md", "", "", "a.md") | Out-Null
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-29 18:40:38"
$wsDe.Range("K3").Value = "2016-08-29 18:40:38"

# ---------------------------------------------------------------------------
# 3. Column widths grow to fit the new, longer text.
# ---------------------------------------------------------------------------
# ColumnWidth is in "characters"; OOXML stores width + 5/6 (~0.8333) more
# than the value assigned here, so we back that padding out to land on the
# target stored widths.
$wsOverview.Columns.Item(5).ColumnWidth = 29.1443713960194   # -> ~29.98
$wsOverview.Columns.Item(6).ColumnWidth = 29.1443713960194   # -> ~29.98

$wsZh.Columns.Item(3).ColumnWidth = 29.1443713960194         # -> ~29.98
$wsZh.Columns.Item(10).ColumnWidth = 39.1666666666667        # -> 40

$wsDe.Columns.Item(3).ColumnWidth = 29.1443713960194         # -> ~29.98
$wsDe.Columns.Item(10).ColumnWidth = 39.1666666666667        # -> 40
